# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# Mirrors the commit: "Created functions to get season record" - the
# existing sheet only had team/player stats, this adds the team's
# win-loss-tie record as three new trailing columns (AD:AF) applied to
# every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing header cell (AC1, style "1":
# bold, centered, thin box border) onto the three new header cells so
# they look like the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record is the same for every player on the roster - 67 wins,
# 95 losses, 0 ties - so stamp it down column AD/AE/AF for every data row.
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 67
    $ws.Cells.Item($r, 31).Value = 95
    $ws.Cells.Item($r, 32).Value = 0
}
